# Rework the "jamming icon" slide: move the wifi picture behind the
# lightning bolt shape (z-order) and re-position/resize both shapes to
# their new layout.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$lightningBolt = $s.Shapes.Item("Lightning Bolt 19")
$picture = $s.Shapes.Item("Picture 24")

# --- Picture 24 (wifi icon): new position/size ---
$picture.Left = 0.0
$picture.Top = 66.0
$picture.Width = 365.01165771484375
$picture.Height = 242.32708740234375

# Send the picture behind the lightning bolt (it used to be drawn after/
# on top of it; now it sits behind it in the z-order).
$picture.ZOrder(1)  # msoSendToBack

# --- Lightning Bolt 19: new rotation/position/size ---
$lightningBolt.Rotation = 44.44649887084961
$lightningBolt.Left = 260.0929260253906
$lightningBolt.Top = 52.77984619140625
$lightningBolt.Width = 361.7464599609375
$lightningBolt.Height = 339.4767761230469
